$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.932.16"
$ws.Range('E2').Value = "'  -3.31%  "
$ws.Range('D3').Value = "'3.358.77"
$ws.Range('E3').Value = "'  -2.82%  "
$ws.Range('E4').Value = "'  -0.03%  "
$ws.Range('D5').Value = "'565.80"
$ws.Range('E5').Value = "'  -2.36%  "
$ws.Range('D6').Value = "'148.24"
$ws.Range('E6').Value = "'  -0.56%  "
$ws.Range('E7').Value = "'  +0.01%  "
$ws.Range('E8').Value = "'  -0.07%  "
$ws.Range('D9').Value = "'7.97"
$ws.Range('E9').Value = "'  +0.88%  "
$ws.Range('E10').Value = "'  -1.58%  "
$ws.Range('E11').Value = "'  +1.29%  "
$ws.Range('D12').Value = "'3.935.16"
$ws.Range('E12').Value = "'  -2.81%  "
$ws.Range('E13').Value = "'  +0.60%  "
$ws.Range('D14').Value = "'27.97"
$ws.Range('E14').Value = "'  -1.77%  "
$ws.Range('D15').Value = "'3.364.04"
$ws.Range('E15').Value = "'  -2.58%  "
$ws.Range('E16').Value = "'  -1.88%  "
$ws.Range('D17').Value = "'61.009.83"
$ws.Range('E17').Value = "'  -3.24%  "
$ws.Range('D18').Value = "'6.33"
$ws.Range('E18').Value = "'  -2.47%  "
$ws.Range('D19').Value = "'14.23"
$ws.Range('E19').Value = "'  -2.96%  "
$ws.Range('E20').Value = "'  -4.21%  "
$ws.Range('D21').Value = "'374.92"
$ws.Range('E21').Value = "'  -3.60%  "
$ws.Range('D22').Value = "'75.27"
$ws.Range('E22').Value = "'  +0.76%  "
$ws.Range('E23').Value = "'  -1.01%  "
$ws.Range('E24').Value = "'  -0.01%  "
$ws.Range('D25').Value = "'3.500.63"
$ws.Range('E25').Value = "'  -2.62%  "
$ws.Range('E26').Value = "'  -6.71%  "
$ws.Range('E27').Value = "'  -3.78%  "
$ws.Range('D28').Value = "'0.997"
$ws.Range('E28').Value = "'  -0.27%  "
$ws.Range('D29').Value = "'7.39"
$ws.Range('E29').Value = "'  -4.14%  "
$ws.Range('E30').Value = "'  -0.07%  "
$ws.Range('E31').Value = "'  -1.80%  "
$ws.Range('E32').Value = "'  -4.65%  "
$ws.Range('D33').Value = "'22.83"
$ws.Range('E33').Value = "'  -2.18%  "
$ws.Range('E34').Value = "'  -4.60%  "
$ws.Range('D35').Value = "'5.36"
$ws.Range('E35').Value = "'  +0.16%  "
$ws.Range('D36').Value = "'170.25"
$ws.Range('E36').Value = "'  +0.04%  "
$ws.Range('E37').Value = "'  -6.04%  "
$ws.Range('E38').Value = "'  -4.09%  "
$ws.Range('D39').Value = "'28.88"
$ws.Range('E39').Value = "'  -10.12%  "
$ws.Range('D40').Value = "'3.393.40"
$ws.Range('E40').Value = "'  -2.85%  "
$ws.Range('E41').Value = "'  -4.60%  "
$ws.Range('D42').Value = "'42.34"
$ws.Range('E42').Value = "'  -1.33%  "
$ws.Range('D43').Value = "'0.759"
$ws.Range('E43').Value = "'  -4.29%  "
$ws.Range('E44').Value = "'  -1.77%  "
$ws.Range('E45').Value = "'  -3.25%  "
$ws.Range('E46').Value = "'  -6.38%  "
$ws.Range('D47').Value = "'2.487.33"
$ws.Range('E47').Value = "'  -2.88%  "
$ws.Range('E48').Value = "'  -3.34%  "
$ws.Range('B49').Value = "'InjectiveProtocol"
$ws.Range('C49').Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range('D49').Value = "'22.44"
$ws.Range('E49').Value = "'  -1.15%  "
$ws.Range('B50').Value = "'FirstDigitalUSD"
$ws.Range('C50').Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range('D50').Value = "'1.00"
$ws.Range('E50').Value = "'  +0.00%  "
$ws.Range('E51').Value = "'  -2.60%  "
